# Weekly update: insert a new price-report row at row 7 (pushing the
# existing rows 7-43 down to 8-44) and populate it with the latest
# "Arveja Verde" market data for Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7:43 shift down to 8:44.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with this week's data.
$ws.Cells.Item(7, 1).Value  = 11
$ws.Cells.Item(7, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value  = "Bíobío"
$ws.Cells.Item(7, 4).Value  = 45092
$ws.Cells.Item(7, 5).Value  = 8
$ws.Cells.Item(7, 6).Value  = 100112022
$ws.Cells.Item(7, 7).Value  = "Arveja Verde"
$ws.Cells.Item(7, 8).Value  = "Sin especificar"
$ws.Cells.Item(7, 9).Value  = "Primera"
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 32000
$ws.Cells.Item(7, 12).Value = 34000
$ws.Cells.Item(7, 13).Value = 33000
$ws.Cells.Item(7, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1320
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
